$d = $word.ActiveDocument
$d.PageSetup.Orientation = 1
